{"js": "// Replace the date heading and all 100 arithmetic answers in the table,\n// matching paragraphs by their position in document order (the date\n// paragraph followed by the 20x5 table of \"<expr>=<result>\" cells).\nconst NEW_VALUES = [\"2024-01-22 Monday\", \"51-48=3\", \"33-20=13\", \"31+26=57\", \"33+50=83\", \"88-64=24\", \"24+26=50\", \"74-22=52\", \"27+67=94\", \"76-72=4\", \"37+29=66\", \"3+46=49\", \"87-0=87\", \"86-43=43\", \"93-30=63\", \"6-3=3\", \"51+38=89\", \"57-2=55\", \"58+36=94\", \"20+22=42\", \"88-39=49\", \"52+35=87\", \"55-24=31\", \"97-70=27\", \"78+9=87\", \"65-26=39\", \"92-75=17\", \"98-84=14\", \"0+66=66\", \"31-8=23\", \"37+19=56\", \"61-60=1\", \"26-0=26\", \"31-7=24\", \"19+64=83\", \"28+2=30\", \"13+78=91\", \"64-47=17\", \"91-19=72\", \"6+5=11\", \"34+44=78\", \"94-94=0\", \"5+0=5\", \"82-68=14\", \"81-29=52\", \"4+28=32\", \"29+1=30\", \"96-52=44\", \"51+25=76\", \"50+5=55\", \"4+27=31\", \"44-20=24\", \"87-73=14\", \"71-27=44\", \"71-36=35\", \"94-71=23\", \"52+23=75\", \"74-73=1\", \"97-88=9\", \"76-1=75\", \"74-12=62\", \"65+20=85\", \"7+58=65\", \"22+30=52\", \"78-40=38\", \"69-7=62\", \"65-10=55\", \"69+0=69\", \"15+24=39\", \"32+66=98\", \"5+33=38\", \"66-2=64\", \"25-24=1\", \"89-28=61\", \"59-9=50\", \"79-40=39\", \"5+64=69\", \"33-3=30\", \"46+39=85\", \"63-19=44\", \"38-20=18\", \"82+15=97\", \"10+81=91\", \"65-25=40\", \"82-59=23\", \"8+22=30\", \"78+14=92\", \"0+68=68\", \"3+30=33\", \"76-6=70\", \"20-4=16\", \"20+0=20\", \"86-24=62\", \"41+54=95\", \"77-1=76\", \"32+33=65\", \"55-39=16\", \"1+7=8\", \"89-89=0\", \"67-62=5\", \"1+53=54\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = Math.min(NEW_VALUES.length, paragraphs.items.length);\nfor (let i = 0; i < count; i++) {\n  paragraphs.items[i].getRange().insertText(NEW_VALUES[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading (first paragraph, outside the table).\n$d.Paragraphs.Item(1).Range.Text = \"2024-01-22 Monday\"\n\n# Update each answer cell in the 20x5 practice table, row by row,\n# matching the before/after order of the source diff.\n$answers = @(\n    @(\"51-48=3\", \"33-20=13\", \"31+26=57\", \"33+50=83\", \"88-64=24\"),\n    @(\"24+26=50\", \"74-22=52\", \"27+67=94\", \"76-72=4\", \"37+29=66\"),\n    @(\"3+46=49\", \"87-0=87\", \"86-43=43\", \"93-30=63\", \"6-3=3\"),\n    @(\"51+38=89\", \"57-2=55\", \"58+36=94\", \"20+22=42\", \"88-39=49\"),\n    @(\"52+35=87\", \"55-24=31\", \"97-70=27\", \"78+9=87\", \"65-26=39\"),\n    @(\"92-75=17\", \"98-84=14\", \"0+66=66\", \"31-8=23\", \"37+19=56\"),\n    @(\"61-60=1\", \"26-0=26\", \"31-7=24\", \"19+64=83\", \"28+2=30\"),\n    @(\"13+78=91\", \"64-47=17\", \"91-19=72\", \"6+5=11\", \"34+44=78\"),\n    @(\"94-94=0\", \"5+0=5\", \"82-68=14\", \"81-29=52\", \"4+28=32\"),\n    @(\"29+1=30\", \"96-52=44\", \"51+25=76\", \"50+5=55\", \"4+27=31\"),\n    @(\"44-20=24\", \"87-73=14\", \"71-27=44\", \"71-36=35\", \"94-71=23\"),\n    @(\"52+23=75\", \"74-73=1\", \"97-88=9\", \"76-1=75\", \"74-12=62\"),\n    @(\"65+20=85\", \"7+58=65\", \"22+30=52\", \"78-40=38\", \"69-7=62\"),\n    @(\"65-10=55\", \"69+0=69\", \"15+24=39\", \"32+66=98\", \"5+33=38\"),\n    @(\"66-2=64\", \"25-24=1\", \"89-28=61\", \"59-9=50\", \"79-40=39\"),\n    @(\"5+64=69\", \"33-3=30\", \"46+39=85\", \"63-19=44\", \"38-20=18\"),\n    @(\"82+15=97\", \"10+81=91\", \"65-25=40\", \"82-59=23\", \"8+22=30\"),\n    @(\"78+14=92\", \"0+68=68\", \"3+30=33\", \"76-6=70\", \"20-4=16\"),\n    @(\"20+0=20\", \"86-24=62\", \"41+54=95\", \"77-1=76\", \"32+33=65\"),\n    @(\"55-39=16\", \"1+7=8\", \"89-89=0\", \"67-62=5\", \"1+53=54\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $answers[$r - 1][$c - 1]\n    }\n}\n"}
